$wb = $excel.ActiveWorkbook

# This script applies numeric updates to several worksheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) reflecting refreshed market-board pricing data pulled by the
# scheduled sheet-update runner. Each block targets one data row identified by its
# leve name (row number) and updates the currentAveragePrice / NQ / HQ / profit columns.

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 93
$ws.Cells.Item(93, 8).Value = 46875
$ws.Cells.Item(93, 10).Value = 46875
$ws.Cells.Item(93, 12).Value = 46875
$ws.Cells.Item(93, 14).Value = -51867
# Row 103
$ws.Cells.Item(103, 8).Value = 1050
$ws.Cells.Item(103, 9).Value = 750
$ws.Cells.Item(103, 11).Value = 2250
$ws.Cells.Item(103, 13).Value = -1664
# Row 106
$ws.Cells.Item(106, 8).Value = 3064.8333
$ws.Cells.Item(106, 9).Value = 2717.8
$ws.Cells.Item(106, 11).Value = 2717.8
$ws.Cells.Item(106, 13).Value = -2086.8
# Row 109
$ws.Cells.Item(109, 8).Value = 60261.332
$ws.Cells.Item(109, 10).Value = 60261.332
$ws.Cells.Item(109, 12).Value = 60261.332
$ws.Cells.Item(109, 14).Value = -63035.332
# Row 128
$ws.Cells.Item(128, 8).Value = 40000
$ws.Cells.Item(128, 10).Value = 40000
$ws.Cells.Item(128, 12).Value = 40000
$ws.Cells.Item(128, 14).Value = -49960
# Row 132
$ws.Cells.Item(132, 8).Value = 4786.6665
$ws.Cells.Item(132, 9).Value = 4786.6665
$ws.Cells.Item(132, 11).Value = 14359.9995
$ws.Cells.Item(132, 13).Value = -11829.9995
# Row 137
$ws.Cells.Item(137, 8).Value = 531690.2
$ws.Cells.Item(137, 9).Value = 2955.6155
$ws.Cells.Item(137, 10).Value = 903233.4399999999
$ws.Cells.Item(137, 11).Value = 8866.8465
$ws.Cells.Item(137, 12).Value = 2709700.32
$ws.Cells.Item(137, 13).Value = -6316.8465
$ws.Cells.Item(137, 14).Value = -2714800.32
# Row 138
$ws.Cells.Item(138, 8).Value = 5289.811
$ws.Cells.Item(138, 9).Value = 2193.6875
$ws.Cells.Item(138, 10).Value = 7648.7617
$ws.Cells.Item(138, 11).Value = 6581.0625
$ws.Cells.Item(138, 12).Value = 22946.2851
$ws.Cells.Item(138, 13).Value = -1441.0625
$ws.Cells.Item(138, 14).Value = -33226.2851

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Cells.Item(74, 8).Value = 6443.4
$ws.Cells.Item(74, 9).Value = 2415.9412
$ws.Cells.Item(74, 10).Value = 15001.75
$ws.Cells.Item(74, 11).Value = 2415.9412
$ws.Cells.Item(74, 12).Value = 15001.75
$ws.Cells.Item(74, 13).Value = -1541.9412
$ws.Cells.Item(74, 14).Value = -16749.75
# Row 77
$ws.Cells.Item(77, 8).Value = 6443.4
$ws.Cells.Item(77, 9).Value = 2415.9412
$ws.Cells.Item(77, 10).Value = 15001.75
$ws.Cells.Item(77, 11).Value = 12079.706
$ws.Cells.Item(77, 12).Value = 75008.75
$ws.Cells.Item(77, 13).Value = -7711.706000000002
$ws.Cells.Item(77, 14).Value = -83744.75
# Row 132
$ws.Cells.Item(132, 8).Value = 2579.0256
$ws.Cells.Item(132, 9).Value = 1628.3
$ws.Cells.Item(132, 10).Value = 3579.7896
$ws.Cells.Item(132, 11).Value = 4884.9
$ws.Cells.Item(132, 12).Value = 10739.3688
$ws.Cells.Item(132, 13).Value = -2354.9
$ws.Cells.Item(132, 14).Value = -15799.3688

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Cells.Item(107, 8).Value = 2136.4375
$ws.Cells.Item(107, 9).Value = 1314.8334
$ws.Cells.Item(107, 10).Value = 2629.4
$ws.Cells.Item(107, 11).Value = 1314.8334
$ws.Cells.Item(107, 12).Value = 2629.4
$ws.Cells.Item(107, 13).Value = 605.1666
$ws.Cells.Item(107, 14).Value = -6469.4

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 621306.5
$ws.Cells.Item(31, 9).Value = 11247.842
$ws.Cells.Item(31, 10).Value = 918514.5600000001
$ws.Cells.Item(31, 11).Value = 11247.842
$ws.Cells.Item(31, 12).Value = 918514.5600000001
$ws.Cells.Item(31, 13).Value = -10952.842
$ws.Cells.Item(31, 14).Value = -919104.5600000001
# Row 34
$ws.Cells.Item(34, 8).Value = 621306.5
$ws.Cells.Item(34, 9).Value = 11247.842
$ws.Cells.Item(34, 10).Value = 918514.5600000001
$ws.Cells.Item(34, 11).Value = 11247.842
$ws.Cells.Item(34, 12).Value = 918514.5600000001
$ws.Cells.Item(34, 13).Value = -11045.842
$ws.Cells.Item(34, 14).Value = -918918.5600000001
# Row 58
$ws.Cells.Item(58, 8).Value = 3032920
$ws.Cells.Item(58, 9).Value = 4547752
$ws.Cells.Item(58, 10).Value = 3255.4
$ws.Cells.Item(58, 11).Value = 4547752
$ws.Cells.Item(58, 12).Value = 3255.4
$ws.Cells.Item(58, 13).Value = -4547549
$ws.Cells.Item(58, 14).Value = -3661.4
# Row 60
$ws.Cells.Item(60, 8).Value = 20966.666
$ws.Cells.Item(60, 10).Value = 20966.666
$ws.Cells.Item(60, 12).Value = 20966.666
$ws.Cells.Item(60, 14).Value = -21988.666
# Row 134
$ws.Cells.Item(134, 8).Value = 2238.775
$ws.Cells.Item(134, 9).Value = 1774.6177
$ws.Cells.Item(134, 11).Value = 5323.8531
$ws.Cells.Item(134, 13).Value = -2788.8531
# Row 136
$ws.Cells.Item(136, 8).Value = 3032920
$ws.Cells.Item(136, 9).Value = 4547752
$ws.Cells.Item(136, 10).Value = 3255.4
$ws.Cells.Item(136, 11).Value = 13643256
$ws.Cells.Item(136, 12).Value = 9766.200000000001
$ws.Cells.Item(136, 13).Value = -13640706
$ws.Cells.Item(136, 14).Value = -14866.2

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Cells.Item(68, 8).Value = 1400.9
$ws.Cells.Item(68, 9).Value = 1454.0435
$ws.Cells.Item(68, 10).Value = 1355.6296
$ws.Cells.Item(68, 11).Value = 4362.1305
$ws.Cells.Item(68, 12).Value = 4066.8888
$ws.Cells.Item(68, 13).Value = -3551.1305
$ws.Cells.Item(68, 14).Value = -5688.8888
# Row 71
$ws.Cells.Item(71, 8).Value = 1400.9
$ws.Cells.Item(71, 9).Value = 1454.0435
$ws.Cells.Item(71, 10).Value = 1355.6296
$ws.Cells.Item(71, 11).Value = 13086.3915
$ws.Cells.Item(71, 12).Value = 12200.6664
$ws.Cells.Item(71, 13).Value = -9030.3915
$ws.Cells.Item(71, 14).Value = -20312.6664
# Row 131
$ws.Cells.Item(131, 8).Value = 39159
$ws.Cells.Item(131, 10).Value = 239225
$ws.Cells.Item(131, 12).Value = 717675
$ws.Cells.Item(131, 14).Value = -727755
# Row 132
$ws.Cells.Item(132, 8).Value = 3643.889
$ws.Cells.Item(132, 9).Value = 3160
$ws.Cells.Item(132, 10).Value = 4248.75
$ws.Cells.Item(132, 11).Value = 28440
$ws.Cells.Item(132, 12).Value = 38238.75
$ws.Cells.Item(132, 13).Value = -25910
$ws.Cells.Item(132, 14).Value = -43298.75

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Cells.Item(107, 8).Value = 7187.1816
$ws.Cells.Item(107, 9).Value = 36150.5
$ws.Cells.Item(107, 10).Value = 750.8889
$ws.Cells.Item(107, 11).Value = 36150.5
$ws.Cells.Item(107, 12).Value = 750.8889
$ws.Cells.Item(107, 13).Value = -34230.5
$ws.Cells.Item(107, 14).Value = -4590.8889

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 74
$ws.Cells.Item(74, 8).Value = 40144.668
$ws.Cells.Item(74, 9).Value = 48000
$ws.Cells.Item(74, 10).Value = 36217
$ws.Cells.Item(74, 11).Value = 48000
$ws.Cells.Item(74, 12).Value = 36217
$ws.Cells.Item(74, 13).Value = -47002
$ws.Cells.Item(74, 14).Value = -38213
# Row 77
$ws.Cells.Item(77, 8).Value = 40144.668
$ws.Cells.Item(77, 9).Value = 48000
$ws.Cells.Item(77, 10).Value = 36217
$ws.Cells.Item(77, 11).Value = 144000
$ws.Cells.Item(77, 12).Value = 108651
$ws.Cells.Item(77, 13).Value = -139008
$ws.Cells.Item(77, 14).Value = -118635
# Row 136
$ws.Cells.Item(136, 8).Value = 5904.5757
$ws.Cells.Item(136, 9).Value = 2446.625
$ws.Cells.Item(136, 10).Value = 9159.117
$ws.Cells.Item(136, 11).Value = 7339.875
$ws.Cells.Item(136, 12).Value = 27477.351
$ws.Cells.Item(136, 13).Value = -4789.875
$ws.Cells.Item(136, 14).Value = -32577.351

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 74
$ws.Cells.Item(74, 8).Value = 15924.777
$ws.Cells.Item(74, 10).Value = 16665.375
$ws.Cells.Item(74, 12).Value = 16665.375
$ws.Cells.Item(74, 14).Value = -18537.375
# Row 75
$ws.Cells.Item(75, 8).Value = 36038
$ws.Cells.Item(75, 10).Value = 36038
$ws.Cells.Item(75, 12).Value = 36038
$ws.Cells.Item(75, 14).Value = -37910
# Row 77
$ws.Cells.Item(77, 8).Value = 15924.777
$ws.Cells.Item(77, 10).Value = 16665.375
$ws.Cells.Item(77, 12).Value = 49996.125
$ws.Cells.Item(77, 14).Value = -59356.125
# Row 78
$ws.Cells.Item(78, 8).Value = 36038
$ws.Cells.Item(78, 10).Value = 36038
$ws.Cells.Item(78, 12).Value = 108114
$ws.Cells.Item(78, 14).Value = -117474
# Row 132
$ws.Cells.Item(132, 8).Value = 5812.4
$ws.Cells.Item(132, 9).Value = 6252.615
$ws.Cells.Item(132, 10).Value = 2951
$ws.Cells.Item(132, 11).Value = 18757.845
$ws.Cells.Item(132, 12).Value = 8853
$ws.Cells.Item(132, 13).Value = -16227.845
$ws.Cells.Item(132, 14).Value = -13913
# Row 136
$ws.Cells.Item(136, 8).Value = 6498.0806
$ws.Cells.Item(136, 9).Value = 6427.0645
$ws.Cells.Item(136, 10).Value = 6569.0967
$ws.Cells.Item(136, 11).Value = 19281.1935
$ws.Cells.Item(136, 12).Value = 19707.2901
$ws.Cells.Item(136, 13).Value = -16731.1935
$ws.Cells.Item(136, 14).Value = -24807.2901
